$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose values change. Excel treats a numeric-looking string assigned
# to Range.Value as a number unless the cell is already formatted as Text,
# so force "@" (text) format first, then restore formatting afterwards —
# this keeps every touched cell stored as a shared string (t="s"), matching
# the rest of row 2, instead of turning it into a numeric cell.
$range = $ws.Range("A2:O2")
$range.NumberFormat = "@"

$ws.Range("A2").Value = "140174"
$ws.Range("D2").Value = "1600"
$ws.Range("E2").Value = "1"
$ws.Range("F2").Value = "1"
$ws.Range("G2").Value = "6"
$ws.Range("H2").Value = "8000"
$ws.Range("J2").Value = "53"
$ws.Range("L2").Value = "0"
$ws.Range("M2").Value = "0"
$ws.Range("O2").Value = "0"

$range.ClearFormats()
